$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Tipo*" column (D) and the legend column (F); table becomes A:C only
$ws.Range("D1:D23").Clear()
$ws.Range("F1:F23").Clear()

# Rewrite the message table (rows 2-7) with the new ordering/codes, and
# append the two new disconnection messages as rows 8 and 9.
$ws.Range("A2").Value = "SendDataConfirmation"
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = "Manda il GUID che è stata generata al client"

$ws.Range("A3").Value = "SendUserNameToHost"
$ws.Range("B3").Value = 1005
$ws.Range("C3").Value = "Manda all'host lo username una volta ricevuta il GUID"

$ws.Range("A4").Value = "LobbyInfoMessage"
$ws.Range("B4").Value = 1010
$ws.Range("C4").Value = "L'Host manda al client le informazioni relative alla chat (HostIp e Lista utenti)"

$ws.Range("A5").Value = "SendUpdatedUserList"
$ws.Range("B5").Value = 1015
$ws.Range("C5").Value = "Invia la lista aggiornata di utenti quanto se ne connette uno nuovo o si disconnette un nuovo utente"

$ws.Range("A6").Value = "LobbyChatMessage"
$ws.Range("B6").Value = 1020
$ws.Range("C6").Value = "Contiene un messaggio testuale per della chat della lobby, e le informazioni del mittente"

$ws.Range("A7").Value = "LobbyStatusAndSettings"
$ws.Range("B7").Value = 1025
$ws.Range("C7").Value = "Contiene lo stato e le impostazioni della partita nella lobby"

$ws.Range("A8").Value = "HostDisconnectedMessage"
$ws.Range("B8").Value = 1026
$ws.Range("C8").Value = "Notifica i client che l'host si è disconnesso"

$ws.Range("A9").Value = "ClientDisconnectedMessage"
$ws.Range("B9").Value = 1027
$ws.Range("C9").Value = "Notifica l'host che il client si è disconnesso manualmente"

[void]$ws.Range("C9").Select()
